$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.795.02'
$ws.Range("E2").Value = '  -0.55%  '

$ws.Range("D3").Value = '3.843.72'
$ws.Range("E3").Value = '  +2.54%  '

$ws.Range("E4").Value = '  -0.09%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '601.10'
$ws.Range("E5").Value = '  -0.15%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '165.99'
$ws.Range("E6").Value = '  -1.14%  '

$ws.Range("D7").Value = '3.840.23'
$ws.Range("E7").Value = '  +2.54%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.530'
$ws.Range("E9").Value = '  -2.01%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.168'
$ws.Range("E10").Value = '  -2.55%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.37'
$ws.Range("E11").Value = '  -0.84%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.460'
$ws.Range("E12").Value = '  -0.38%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '36.95'
$ws.Range("E13").Value = '  -3.10%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000244'
$ws.Range("E14").Value = '  -2.02%  '

$ws.Range("D15").Value = '4.487.13'
$ws.Range("E15").Value = '  +2.54%  '

$ws.Range("D16").Value = '3.854.98'
$ws.Range("E16").Value = '  +3.06%  '

$ws.Range("D17").Value = '68.852.85'
$ws.Range("E17").Value = '  -0.53%  '

$ws.Range("E18").Value = '  +0.81%  '

$ws.Range("E19").Value = '  -0.64%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.09'
$ws.Range("E20").Value = '  -1.69%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.26'
$ws.Range("E21").Value = '  -0.60%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '487.08'
$ws.Range("E22").Value = '  -1.25%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.718'
$ws.Range("E23").Value = '  -1.91%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.0000165'
$ws.Range("E24").Value = '  +8.19%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '84.07'
$ws.Range("E25").Value = '  -1.05%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.25'
$ws.Range("E26").Value = '  -2.18%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.10'
$ws.Range("E27").Value = '  -1.97%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.00'
$ws.Range("E28").Value = '  +0.04%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '10.00'
$ws.Range("E29").Value = '  -0.65%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.95'
$ws.Range("E30").Value = '  -1.35%  '

$ws.Range("D31").Value = '3.993.97'
$ws.Range("E31").Value = '  +2.56%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.88'
$ws.Range("E32").Value = '  -4.94%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.37'
$ws.Range("E33").Value = '  -4.54%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '31.73'
$ws.Range("E34").Value = '  +0.39%  '

$ws.Range("D35").Value = '3.789.30'
$ws.Range("E35").Value = '  +2.92%  '

$ws.Range("E36").Value = '  -1.27%  '

$ws.Range("E37").Value = '  +1.65%  '

$ws.Range("E38").Value = '  -0.31%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.88'
$ws.Range("E39").Value = '  -2.01%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.00'
$ws.Range("E40").Value = '  -0.04%  '

$ws.Range("B41").Value = 'Bittensor'
$ws.Range("C41").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '441.89'
$ws.Range("E41").Value = '  +3.92%  '

$ws.Range("B42").Value = 'TheGraph'
$ws.Range("C42").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.318'
$ws.Range("E42").Value = '  -2.74%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.02'
$ws.Range("E43").Value = '  -4.94%  '

$ws.Range("B44").Value = 'Stacks'
$ws.Range("C44").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.00'
$ws.Range("E44").Value = '  +0.10%  '

$ws.Range("B45").Value = 'OKB'
$ws.Range("C45").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '48.54'
$ws.Range("E45").Value = '  -0.48%  '

$ws.Range("E46").Value = '  +0.02%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.44'
$ws.Range("E47").Value = '  -0.64%  '

$ws.Range("B48").Value = 'Monero'
$ws.Range("C48").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '142.09'
$ws.Range("E48").Value = '  +0.39%  '

$ws.Range("B49").Value = 'Maker'
$ws.Range("C49").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D49").Value = '2.815.20'
$ws.Range("E49").Value = '  +0.68%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0358'
$ws.Range("E50").Value = '  +0.75%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '25.91'
$ws.Range("E51").Value = '  +11.60%  '
